$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New data rows appended at the bottom of the sheet (rows 100 and 101),
# mirroring the formatting used by the existing data rows (e.g. row 99).
# ---------------------------------------------------------------------------

$newRows = @(
    @{
        Row = 100
        Values = @{
            A = 99
            B = "ecuador"
            C = "liga-pro"
            D = "2023"
            E = 45241.89583333334
            F = "Guayaquil City"
            G = 1
            H = "Tecnico U."
            I = 0
            J = 3.78
            K = "09/11/2023 01:12"
            L = 3.29
            M = "11/11/2023 21:29"
            N = 3.28
            O = "09/11/2023 01:12"
            P = 3.14
            Q = "11/11/2023 21:26"
            R = 1.99
            S = "09/11/2023 01:12"
            T = 2.37
            U = "11/11/2023 21:26"
            V = "https://www.betexplorer.com/football/ecuador/liga-pro/guayaquil-city-tecnico-u/MDDvFw4q/"
        }
    },
    @{
        Row = 101
        Values = @{
            A = 100
            B = "ecuador"
            C = "liga-pro"
            D = "2023"
            E = 45242
            F = "Aucas"
            G = 2
            H = "Ind. del Valle"
            I = 0
            J = 2.94
            K = "05/11/2023 00:12"
            L = 3.21
            M = "11/11/2023 23:53"
            N = 3.28
            O = "05/11/2023 00:12"
            P = 3.52
            Q = "11/11/2023 23:53"
            R = 2.44
            S = "05/11/2023 00:12"
            T = 2.23
            U = "11/11/2023 23:53"
            V = "https://www.betexplorer.com/football/ecuador/liga-pro/aucas-independiente-del-valle/2P06XoKA/"
        }
    }
)

# Source row whose per-column formatting (bold/border on "Indice", date/time
# number format on "data_partida") will be copied onto the new rows.
$templateRow = 99

foreach ($entry in $newRows) {
    $r = $entry.Row
    $vals = $entry.Values

    # Copy the formatting of column A (Indice: bold, centered, bordered) and
    # column E (data_partida: date/time number format) from the template row.
    $ws.Range("A$templateRow").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)

    $ws.Range("E$templateRow").Copy()
    $ws.Range("E$r").PasteSpecial(-4122)

    $ws.Range("A$r").Value = $vals.A
    $ws.Range("E$r").Value = $vals.E

    $ws.Range("B$r").Value = $vals.B
    $ws.Range("C$r").Value = $vals.C

    # Column D holds a purely numeric-looking string ("2023") that must stay
    # text, like the rest of the column; force text entry then drop the
    # temporary formatting so no new style is left behind.
    $ws.Range("D$r").NumberFormat = "@"
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("D$r").ClearFormats()

    $ws.Range("F$r").Value = $vals.F
    $ws.Range("G$r").Value = $vals.G
    $ws.Range("H$r").Value = $vals.H
    $ws.Range("I$r").Value = $vals.I
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("Q$r").Value = $vals.Q
    $ws.Range("R$r").Value = $vals.R
    $ws.Range("S$r").Value = $vals.S
    $ws.Range("T$r").Value = $vals.T
    $ws.Range("U$r").Value = $vals.U
    $ws.Range("V$r").Value = $vals.V
}

Write-Host "Added rows 100-101"
